$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Three tables (slides 14, 15, 16) switch from the deck's custom
#    "Table_0" style ({7A4DE897-BD1A-4DB0-8F8B-1D1FC60199DD}) to the
#    built-in table style {CBCF4722-E73B-4F44-AEFC-4AFC234D9E35}.
# ---------------------------------------------------------------------------
$tableSlideIndexes = @(14, 15, 16)
$newTableStyleId = "{CBCF4722-E73B-4F44-AEFC-4AFC234D9E35}"

foreach ($slideIdx in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2) The deck's theme (slide master) switches from the "Integral" / "Red
#    Violet" colour scheme to the standard Office colour scheme.
# ---------------------------------------------------------------------------
$cs = $p.SlideMaster.ColorScheme
$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72
